$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-34 & 37-51: refresh Price (D) and Volume(1h) (E) text columns with the
# latest scraped figures. The Price column holds numeric-looking text (e.g.
# "85.00", "0.999", "39.872.85" as a thousands-grouped string) so a leading
# apostrophe is used to force text entry - matching how Excel itself keeps a
# quote-prefixed value as a string instead of re-parsing it as a Double and
# silently dropping trailing zeros / "." thousands separators.
$ws.Range("D2").Value = "'39.872.85"
$ws.Range("E2").Value = "  -0.54%  "
$ws.Range("D3").Value = "'2.231.69"
$ws.Range("E3").Value = "  -4.93%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'293.56"
$ws.Range("E5").Value = "  -5.70%  "
$ws.Range("D6").Value = "'85.00"
$ws.Range("E6").Value = "  -0.27%  "
$ws.Range("D7").Value = "'0.513"
$ws.Range("E7").Value = "  -2.73%  "
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("D9").Value = "'0.468"
$ws.Range("E9").Value = "  -2.98%  "
$ws.Range("D10").Value = "'0.0797"
$ws.Range("E10").Value = "  -1.34%  "
$ws.Range("D11").Value = "'30.13"
$ws.Range("E11").Value = "  -0.05%  "
$ws.Range("D12").Value = "'48.08"
$ws.Range("E12").Value = "  -8.44%  "
$ws.Range("E13").Value = "  -2.29%  "
$ws.Range("D14").Value = "'6.35"
$ws.Range("E14").Value = "  -1.03%  "
$ws.Range("D15").Value = "'2.576.01"
$ws.Range("E15").Value = "  -4.95%  "
$ws.Range("D16").Value = "'14.16"
$ws.Range("E16").Value = "  -4.58%  "
$ws.Range("D17").Value = "'2.224.56"
$ws.Range("E17").Value = "  -5.97%  "
$ws.Range("D18").Value = "'0.723"
$ws.Range("E18").Value = "  -5.16%  "
$ws.Range("D19").Value = "'39.761.94"
$ws.Range("E19").Value = "  -0.92%  "
$ws.Range("D20").Value = "'0.0₃0890"
$ws.Range("E20").Value = "  -1.23%  "
$ws.Range("D21").Value = "'5.78"
$ws.Range("E21").Value = "  -5.11%  "
$ws.Range("D22").Value = "'65.35"
$ws.Range("E22").Value = "  -4.31%  "
$ws.Range("D23").Value = "'10.51"
$ws.Range("E23").Value = "  -1.23%  "
$ws.Range("D24").Value = "'232.11"
$ws.Range("E24").Value = "  -1.42%  "
$ws.Range("E25").Value = "  +0.05%  "
$ws.Range("D26").Value = "'2.41"
$ws.Range("E26").Value = "  -5.51%  "
$ws.Range("D27").Value = "'1.84"
$ws.Range("E27").Value = "  +0.96%  "
$ws.Range("D28").Value = "'22.92"
$ws.Range("E28").Value = "  -3.18%  "
$ws.Range("E29").Value = "  -0.21%  "
$ws.Range("D30").Value = "'9.22"
$ws.Range("E30").Value = "  -0.53%  "
$ws.Range("D31").Value = "'154.57"
$ws.Range("E31").Value = "  +0.38%  "
$ws.Range("D32").Value = "'32.84"
$ws.Range("E32").Value = "  -5.73%  "
$ws.Range("D33").Value = "'0.999"
$ws.Range("E33").Value = "  -0.29%  "
$ws.Range("D34").Value = "'4.84"
$ws.Range("E34").Value = "  -5.28%  "
# Rows 35/36: ranking reshuffled this run - WEMIXToken and Hedera swapped places
# - together with their refreshed Price/Volume figures.
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").Value = "'0.0706"
$ws.Range("E35").Value = "  -1.74%  "

$ws.Range("B36").Value = "WEMIXToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D36").Value = "'2.36"
$ws.Range("E36").Value = "  -5.08%  "
$ws.Range("D37").Value = "'16.35"
$ws.Range("E37").Value = "  +4.16%  "
$ws.Range("D38").Value = "'0.112"
$ws.Range("E38").Value = "  -1.82%  "
$ws.Range("D39").Value = "'0.0981"
$ws.Range("E39").Value = "  -1.27%  "
$ws.Range("D40").Value = "'2.67"
$ws.Range("E40").Value = "  -4.62%  "
$ws.Range("D41").Value = "'1.66"
$ws.Range("E41").Value = "  -3.69%  "
$ws.Range("D42").Value = "'3.73"
$ws.Range("E42").Value = "  -3.09%  "
$ws.Range("D43").Value = "'1.951.46"
$ws.Range("E43").Value = "  -0.91%  "
$ws.Range("D44").Value = "'2.18"
$ws.Range("E44").Value = "  -2.78%  "
$ws.Range("D45").Value = "'0.0270"
$ws.Range("E45").Value = "  +1.47%  "
$ws.Range("D46").Value = "'9.34"
$ws.Range("E46").Value = "  -2.04%  "
$ws.Range("D47").Value = "'16.26"
$ws.Range("E47").Value = "  -7.94%  "
$ws.Range("D48").Value = "'2.60"
$ws.Range("E48").Value = "  -3.38%  "
$ws.Range("D49").Value = "'2.441.60"
$ws.Range("E49").Value = "  -4.98%  "
$ws.Range("D50").Value = "'70.83"
$ws.Range("E50").Value = "  +0.66%  "
$ws.Range("D51").Value = "'88.87"
$ws.Range("E51").Value = "  -4.64%  "
